$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 138
$ws1.Range("F3").Value = 215
$ws1.Range("F4").Value = 3646
$ws1.Range("F5").Value = 378
$ws1.Range("F7").Value = 433

# Sheet "全部类型" (fourth sheet) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 138
$ws4.Range("F3").Value = 215
$ws4.Range("F4").Value = 3646
$ws4.Range("F5").Value = 378
$ws4.Range("F9").Value = 433
